# Añadido cabecera a Reunion2
# Insert a new underlined heading paragraph ("Daniel Ojeda Velasco")
# at the very top of the document, before the existing first paragraph.

$d = $word.ActiveDocument

# Insert a brand-new paragraph immediately before the current first
# paragraph. Word's InsertParagraphBefore() carries the adjacent
# paragraph's (and run's) formatting onto the freshly minted paragraph
# mark, which here is the single-underline formatting already used by
# the document's other section headers.
$firstPara = $d.Paragraphs(1).Range
$firstPara.InsertParagraphBefore()

# The just-inserted paragraph is now Paragraphs(1); set its text.
$heading = $d.Paragraphs(1).Range
$heading.Text = "Daniel Ojeda Velasco"
